$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in A9 (new shared string "casd")
$ws.Range("A9").Value = "casd"

# Update the selection to A9 (matches diff: activeCell="A9" sqref="A9")
$ws.Range("A9").Select()
